# Update "想去人数" (column F) figures on both the "展览" sheet and the
# "全部类型" aggregate sheet to the newly scraped counts.
# (sheet "全部类型" mirrors "展览" rows shifted down by one row because it
# also includes a row from the "演出" sheet above it.)

$wb = $excel.ActiveWorkbook

$updates = @(
    @{Row = 3;  Value = 736}
    @{Row = 4;  Value = 1447}
    @{Row = 5;  Value = 221}
    @{Row = 7;  Value = 132}
    @{Row = 8;  Value = 6140}
    @{Row = 9;  Value = 68}
    @{Row = 10; Value = 399}
    @{Row = 11; Value = 111}
    @{Row = 12; Value = 5009}
    @{Row = 14; Value = 175}
    @{Row = 15; Value = 1170}
    @{Row = 16; Value = 52}
    @{Row = 17; Value = 356}
    @{Row = 20; Value = 285}
    @{Row = 22; Value = 3504}
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.Value
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($u in $updates) {
    $ws4.Cells.Item($u.Row + 1, 6).Value = $u.Value
}
